$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": drop the "Handoff transform failed" row for the old
# b74e7a9b file, and rename the c36e646e file to the newly generated
# 5e25e143 file (fresh handoff report).
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Clear existing hyperlinks so we can rebuild a clean, consistent set after
# the row shift below (stale hyperlink refs are not auto-fixed by Delete).
$wsOverview.Hyperlinks.Delete()

# Remove the middle row (old "b74e7a9b...md" / "Handoff transform failed"),
# shifting ".localization-config" up from row 4 to row 3.
$wsOverview.Range("A3:C3").Delete(-4162)

# Rename the successful handoff file in row 2.
$wsOverview.Range("A2").Value = "5e25e143-e598-41cf-a632-591b6707b2d9.md"

# Rebuild hyperlinks to match the new layout.
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e87275551b8766a1624a2dfc090c95596f152797/e2e/5e25e143-e598-41cf-a632-591b6707b2d9.md", $null, $null, "5e25e143-e598-41cf-a632-591b6707b2d9.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e87275551b8766a1624a2dfc090c95596f152797/.localization-config", $null, $null, ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "zh-cn": same row removal + rename, plus refresh the handoff
# transform file name and handoff datetime for the regenerated report.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Delete()
$wsZh.Range("A3:I3").Delete(-4162)

$wsZh.Range("A2").Value = "5e25e143-e598-41cf-a632-591b6707b2d9.md"
$wsZh.Range("C2").Value = "5e25e143-e598-41cf-a632-591b6707b2d9.a8d6675a054f762c55cebf59e56572f898520da6.zh-cn.xlf"
$wsZh.Range("D2").Value = "2016-02-15 04:06:29"

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e87275551b8766a1624a2dfc090c95596f152797/e2e/5e25e143-e598-41cf-a632-591b6707b2d9.md", $null, $null, "5e25e143-e598-41cf-a632-591b6707b2d9.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b5f905f68002212b470df7aaf2b1e73b9a104d8c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/5e25e143-e598-41cf-a632-591b6707b2d9.a8d6675a054f762c55cebf59e56572f898520da6.zh-cn.xlf", $null, $null, "5e25e143-e598-41cf-a632-591b6707b2d9.a8d6675a054f762c55cebf59e56572f898520da6.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e87275551b8766a1624a2dfc090c95596f152797/.localization-config", $null, $null, ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "de-de": same treatment, with the de-de handoff timestamp.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Delete()
$wsDe.Range("A3:I3").Delete(-4162)

$wsDe.Range("A2").Value = "5e25e143-e598-41cf-a632-591b6707b2d9.md"
$wsDe.Range("C2").Value = "5e25e143-e598-41cf-a632-591b6707b2d9.a8d6675a054f762c55cebf59e56572f898520da6.de-de.xlf"
$wsDe.Range("D2").Value = "2016-02-15 04:06:43"

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e87275551b8766a1624a2dfc090c95596f152797/e2e/5e25e143-e598-41cf-a632-591b6707b2d9.md", $null, $null, "5e25e143-e598-41cf-a632-591b6707b2d9.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5e534fb73e349c573c7f0963cb407b84a29e1aae/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/5e25e143-e598-41cf-a632-591b6707b2d9.a8d6675a054f762c55cebf59e56572f898520da6.de-de.xlf", $null, $null, "5e25e143-e598-41cf-a632-591b6707b2d9.a8d6675a054f762c55cebf59e56572f898520da6.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e87275551b8766a1624a2dfc090c95596f152797/.localization-config", $null, $null, ".localization-config")
